# Swap the values of columns A, B, D, E, F, G, H between row 3 and row 4
# (columns C, I, J, K, ... are identical between the two rows, so they are left untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$cols = @("A", "B", "D", "E", "F", "G", "H")

foreach ($col in $cols) {
    $cellRow3 = $ws.Range($col + "3")
    $cellRow4 = $ws.Range($col + "4")

    $valRow3 = $cellRow3.Value2
    $valRow4 = $cellRow4.Value2

    $cellRow3.Value2 = $valRow4
    $cellRow4.Value2 = $valRow3
}
